$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Cell values - header + 4 new data rows (TC1/TC2 rows replaced by the new
#    TC3_twitLoginChrome / DatadrivenTest rows, plus a new AppURL column).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "TestCaseName"
$ws.Range("B1").Value = "RunMode"
$ws.Range("C1").Value = "Browser"
$ws.Range("D1").Value = "AppURL"
$ws.Range("E1").Value = "UserName"
$ws.Range("F1").Value = "Password"

$ws.Range("A2").Value = "TC3_twitLoginChrome"
$ws.Range("B2").Value = "Y"
$ws.Range("C2").Value = "Chrome"
$ws.Range("D2").Value = "https://twitter.com/login?lang=en"
$ws.Range("E2").Value = "s1"
$ws.Range("F2").Value = "p1"

$ws.Range("A3").Value = "TC3_twitLoginChrome"
$ws.Range("B3").Value = "Y"
$ws.Range("C3").Value = "Chrome"
$ws.Range("D3").Value = "https://twitter.com/login?lang=en"
$ws.Range("E3").Value = "s2"
$ws.Range("F3").Value = "p2"

$ws.Range("A4").Value = "DatadrivenTest"
$ws.Range("B4").Value = "Y"
$ws.Range("C4").Value = "IE"
$ws.Range("D4").Value = "https://twitter.com/login?lang=en"
$ws.Range("E4").Value = "s3"
$ws.Range("F4").Value = "p3"

$ws.Range("A5").Value = "DatadrivenTest"
$ws.Range("B5").Value = "Y"
$ws.Range("C5").Value = "Chrome"
$ws.Range("D5").Value = "https://twitter.com/login?lang=en"
$ws.Range("E5").Value = "s3"
$ws.Range("F5").Value = "p3"

# ---------------------------------------------------------------------------
# 2. Formatting.
#    Build the cellXfs in the same order the workbook ends up with:
#      s=1 center/default font  (columns B:F body)
#      s=2 bold+center/bold font (header row)
#      s=3 left/default font    (column A body)
#      s=4 center/hyperlink font (column D body, "looks like" a link)
# ---------------------------------------------------------------------------

# s=1 : center alignment across columns B:F (this also temporarily covers the
# header row - it gets upgraded to bold+center below).
$ws.Range("B1:F1048576").HorizontalAlignment = -4108

# s=2 : header row -> bold + center.
$ws.Range("A1:F1").HorizontalAlignment = -4108
$ws.Range("A1:F1").Font.Bold = $true

# s=3 : column A body rows -> left aligned.
$ws.Range("A2:A1048576").HorizontalAlignment = -4131

# s=4 : column D body rows -> hyperlink-styled (center + hyperlink font),
# without leaving a real hyperlink object/relationship behind. Hyperlinks.Add
# only stamps the style onto its own anchor cell, so add one per cell.
$ws.Range("D2:D5").HorizontalAlignment = -4108
$ws.Hyperlinks.Add($ws.Range("D2"), "https://twitter.com/login?lang=en") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://twitter.com/login?lang=en") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://twitter.com/login?lang=en") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "https://twitter.com/login?lang=en") | Out-Null
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 3. Column widths (best effort; engine quantizes to 1/6-character steps).
# ---------------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 19.833333333333336
$ws.Columns("B:B").ColumnWidth = 8.833333333333332
$ws.Columns("C:C").ColumnWidth = 7.500000000000001
$ws.Columns("D:D").ColumnWidth = 31.500000000000004
$ws.Columns("E:E").ColumnWidth = 9.666666666666666
$ws.Columns("F:F").ColumnWidth = 8.666666666666666

# ---------------------------------------------------------------------------
# 4. Selection / view state.
# ---------------------------------------------------------------------------
$ws.Range("E5:F5").Select()

# ---------------------------------------------------------------------------
# 5. Page setup (best effort).
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
